$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.225
$ws.Range("C2").Value = 0.496875
$ws.Range("J2").Value = 0.015625
$ws.Range("P2").Value = 0.153125
$ws.Range("S2").Value = 0.109375

# Row 3
$ws.Range("B3").Value = 0.00558659217877095
$ws.Range("C3").Value = 0.03910614525139665
$ws.Range("J3").Value = 0.05027932960893855
$ws.Range("P3").Value = 0.7541899441340782
$ws.Range("S3").Value = 0.1508379888268156

# Row 4
$ws.Range("J4").Value = 0.08695652173913043
$ws.Range("P4").Value = 0.5652173913043478
$ws.Range("S4").Value = 0.3478260869565217

# Row 5
$ws.Range("J5").Value = 0.25
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.25

# Row 6
$ws.Range("B6").Value = 0.05855855855855856
$ws.Range("D6").Value = 0.009009009009009009
$ws.Range("F6").Value = 0.04954954954954955
$ws.Range("J6").Value = 0.3063063063063063
$ws.Range("O6").Value = 0.02702702702702703
$ws.Range("Q6").Value = 0.1486486486486487
$ws.Range("R6").Value = 0.08108108108108109
$ws.Range("S6").Value = 0.3198198198198198

# Row 7
$ws.Range("B7").Value = 0.126984126984127
$ws.Range("D7").Value = 0.02645502645502645
$ws.Range("F7").Value = 0.0582010582010582
$ws.Range("J7").Value = 0.1216931216931217
$ws.Range("O7").Value = 0.01058201058201058
$ws.Range("Q7").Value = 0.1428571428571428
$ws.Range("R7").Value = 0.0582010582010582
$ws.Range("S7").Value = 0.455026455026455

# Row 8
$ws.Range("B8").Value = 0.1066945606694561
$ws.Range("D8").Value = 0.01882845188284519
$ws.Range("E8").Value = 0.002092050209205021
$ws.Range("F8").Value = 0.06485355648535565
$ws.Range("J8").Value = 0.1108786610878661
$ws.Range("O8").Value = 0.01464435146443515
$ws.Range("Q8").Value = 0.1610878661087866
$ws.Range("R8").Value = 0.08368200836820083
$ws.Range("S8").Value = 0.4372384937238494

# Row 9
$ws.Range("B9").Value = 0.06716417910447761
$ws.Range("D9").Value = 0.007462686567164179
$ws.Range("F9").Value = 0.08955223880597014
$ws.Range("J9").Value = 0.1641791044776119
$ws.Range("O9").Value = 0.03731343283582089
$ws.Range("Q9").Value = 0.1268656716417911
$ws.Range("R9").Value = 0.09701492537313433
$ws.Range("S9").Value = 0.4104477611940299

# Row 10
$ws.Range("B10").Value = 0.1157222665602554
$ws.Range("D10").Value = 0.02314445331205108
$ws.Range("E10").Value = 0.002394253790901836
$ws.Range("F10").Value = 0.07262569832402235
$ws.Range("J10").Value = 0.1332801276935355
$ws.Range("O10").Value = 0.02314445331205108
$ws.Range("Q10").Value = 0.2162809257781325
$ws.Range("R10").Value = 0.07901037509976057
$ws.Range("S10").Value = 0.3343974461292897

# Row 11
$ws.Range("F11").Value = 0.0078125
$ws.Range("G11").Value = 0.14453125
$ws.Range("J11").Value = 0.078125
$ws.Range("K11").Value = 0.1796875
$ws.Range("L11").Value = 0.578125
$ws.Range("S11").Value = 0.01171875

# Row 12
$ws.Range("G12").Value = 0.75
$ws.Range("J12").Value = 0.1730769230769231
$ws.Range("K12").Value = 0.00641025641025641
$ws.Range("L12").Value = 0.05128205128205128
$ws.Range("S12").Value = 0.01923076923076923

# Row 13
$ws.Range("G13").Value = 0.7258064516129032
$ws.Range("J13").Value = 0.2096774193548387
$ws.Range("S13").Value = 0.06451612903225806

# Row 15
$ws.Range("F15").Value = 0.01941747572815534
$ws.Range("H15").Value = 0.1553398058252427
$ws.Range("I15").Value = 0.06310679611650485
$ws.Range("J15").Value = 0.3300970873786408
$ws.Range("K15").Value = 0.05825242718446602
$ws.Range("M15").Value = 0.009708737864077669
$ws.Range("O15").Value = 0.04368932038834952
$ws.Range("S15").Value = 0.3203883495145631

# Row 16
$ws.Range("F16").Value = 0.02463054187192118
$ws.Range("H16").Value = 0.1428571428571428
$ws.Range("I16").Value = 0.04433497536945813
$ws.Range("J16").Value = 0.4187192118226601
$ws.Range("K16").Value = 0.1527093596059113
$ws.Range("M16").Value = 0.01970443349753695
$ws.Range("O16").Value = 0.04926108374384237
$ws.Range("S16").Value = 0.1477832512315271

# Row 17
$ws.Range("F17").Value = 0.01157407407407407
$ws.Range("H17").Value = 0.1805555555555556
$ws.Range("I17").Value = 0.06712962962962964
$ws.Range("J17").Value = 0.4166666666666667
$ws.Range("K17").Value = 0.09259259259259259
$ws.Range("M17").Value = 0.03472222222222222
$ws.Range("N17").Value = 0.002314814814814815
$ws.Range("O17").Value = 0.07175925925925926
$ws.Range("S17").Value = 0.1226851851851852

# Row 18
$ws.Range("F18").Value = 0.005376344086021506
$ws.Range("H18").Value = 0.1720430107526882
$ws.Range("I18").Value = 0.05376344086021505
$ws.Range("J18").Value = 0.4838709677419355
$ws.Range("K18").Value = 0.09139784946236559
$ws.Range("M18").Value = 0.02150537634408602
$ws.Range("N18").Value = 0.005376344086021506
$ws.Range("O18").Value = 0.03763440860215054
$ws.Range("S18").Value = 0.1290322580645161

# Row 19
$ws.Range("F19").Value = 0.01483924154987634
$ws.Range("H19").Value = 0.2547403132728772
$ws.Range("I19").Value = 0.06100577081615829
$ws.Range("J19").Value = 0.361088211046991
$ws.Range("K19").Value = 0.09315746084089035
$ws.Range("M19").Value = 0.03050288540807914
$ws.Range("N19").Value = 0.0008244023083264633
$ws.Range("O19").Value = 0.07089859851607584
$ws.Range("S19").Value = 0.1129431162407255
